# LOM3224.xlsx — rewrite the "Programa resumido / Programa / Avaliacao" block
# (rows 13-24) to the condensed content, drop the old row 25, and fix up the
# row heights to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: "Programa resumido:" now holds "Semestral" -------------------
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# --- Row 14: "Short syllabus:" (label only) --------------------------------
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# --- Row 15: "Programa:" -----------------------------------------------
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"

# --- Row 16: "Syllabus:" (label only) --------------------------------------
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()

# --- Row 17: "Avaliação:" (label only) --------------------------------------
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

# --- Row 18: "Método:" now holds the professor's info ----------------------
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"

# --- Row 19: "Critério:" now holds the old "Método" text --------------------
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("C19").Value = "Aulas expositivas, seminários e exercícios comentados."

# --- Row 20: "Norma de recuperação:" now holds the old "Critério" text ------
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("C20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."

# --- Row 21: "Bibliografia:" now holds the old "Norma de recuperação" text --
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

# --- Row 22: "Requisitos:" (label only, bibliography text dropped) ----------
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

# --- Row 23 / 24: requirement list shifts up one row ------------------------
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOM3231 -  Métodos Experimentais da Física IV  (Requisito)`n"
$ws.Range("C23").Value = "LOM3231 -  Métodos Experimentais da Física IV  (Requisito)`n"

$ws.Range("B24").Value = "LOM3234 -  Óptica Física  (Requisito)`n"
$ws.Range("C24").Value = "LOM3234 -  Óptica Física  (Requisito)`n"

# --- Drop the now-unused last row (old row 25) ------------------------------
$ws.Rows(25).Delete()

# --- Fix up row heights to match the new layout -----------------------------
$ws.Rows(13).RowHeight = 60
$ws.Rows(14).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(16).RowHeight = 120
$ws.Rows(17).AutoFit()
$ws.Rows(18).RowHeight = 60
$ws.Rows(19).RowHeight = 60
$ws.Rows(20).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(22).AutoFit()
$ws.Rows(23).RowHeight = 30
$ws.Rows(24).RowHeight = 30
